$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.554.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.341.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "664.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.47%  "

$ws.Range("E7").Value = "  +12.09%  "

$ws.Range("E8").Value = "  +16.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +27.46%  "

$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.337.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.94%  "

$ws.Range("E12").Value = "  +5.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +19.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000267"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "99.487.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.968.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.89%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.337.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +23.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.13%  "

$ws.Range("E21").Value = "  +2.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "528.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000211"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.440"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +57.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.519.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("E30").Value = "  +17.15%  "

$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +16.05%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.538"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +17.35%  "

$ws.Range("E37").Value = "  +7.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.160"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "529.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.07%  "

$ws.Range("E41").Value = "  -0.79%  "

$ws.Range("E42").Value = "  +5.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0442"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +34.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.824"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.60%  "

$ws.Range("E46").Value = "  +3.02%  "

$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.83%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +22.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "164.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.98%  "

